# The drawing in this document (a <v:group> of VML shapes nested inside a
# <w:pict> run) is opaque to the high-level Word object model exposed by this
# runtime: Shapes.Count / InlineShapes.Count report 0 for it, and the hosting
# paragraph's Range.Text is empty, so Find.Execute / Shapes(...) cannot reach
# the nested <v:shapetype>/<v:shape> markup to edit it directly.
#
# The one lever that *does* reach into that markup is Range.InsertXML - but it
# REPLACES the full contents of the Range it is called on. So we rebuild the
# hosting paragraph's content verbatim and splice in the new connector shape
# (plus the <v:shapetype id="_x0000_t34"> element it points at) in the same
# spot the upstream diff adds it: right after the last <v:rect> ("App-V Agent")
# and right before the closing </v:group>.

$d = $word.ActiveDocument

# This fixture keeps the VML drawing in the document's 2nd paragraph; sanity
# check that assumption against the package XML before touching anything.
$full = $d.WordOpenXML
if (-not $full.Contains("_x0000_s1049")) {
    throw "Expected VML shape '_x0000_s1049' (App-V Agent rect) not found in document."
}
if ($d.Paragraphs.Count -lt 2) {
    throw "Expected at least 2 paragraphs; found " + $d.Paragraphs.Count
}

$target = $d.Paragraphs.Item(2).Range

$xml = @'
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:v="urn:schemas-microsoft-com:vml" xmlns:o="urn:schemas-microsoft-com:office:office" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" w:rsidR="00DC1D4E" w:rsidRDefault="00DC1D4E" w:rsidP="00DC1D4E">
      <w:pPr>
        <w:pStyle w:val="BodyTextIndent"/>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Arial" w:hAnsi="Arial"/>
          <w:b/>
          <w:noProof/>
        </w:rPr>
      </w:r>
      <w:r>
        <w:pict>
          <v:group id="_x0000_s1026" style="width:342pt;height:180.65pt;mso-position-horizontal-relative:char;mso-position-vertical-relative:line" coordorigin="2785,-605" coordsize="6514,3468">
            <v:rect id="_x0000_s1028" style="position:absolute;left:2785;top:-260;width:2228;height:3123" fillcolor="yellow"/>
            <v:rect id="_x0000_s1033" style="position:absolute;left:2859;top:2474;width:2094;height:304">
              <v:textbox>
                <w:txbxContent>
                  <w:p w:rsidR="00DC1D4E" w:rsidRPr="00A524AD" w:rsidRDefault="00DC1D4E" w:rsidP="00DC1D4E">
                    <w:pPr>
                      <w:jc w:val="center"/>
                      <w:rPr>
                        <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
                        <w:sz w:val="16"/>
                        <w:szCs w:val="16"/>
                      </w:rPr>
                    </w:pPr>
                    <w:r w:rsidRPr="00A524AD">
                      <w:rPr>
                        <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
                        <w:sz w:val="16"/>
                        <w:szCs w:val="16"/>
                      </w:rPr>
                      <w:t>Hardware</w:t>
                    </w:r>
                  </w:p>
                </w:txbxContent>
              </v:textbox>
            </v:rect>
            <v:shapetype id="_x0000_t34" coordsize="21600,21600" o:spt="34" o:oned="t" adj="10800" path="m,l@0,0@0,21600,21600,21600e" filled="f">
              <v:stroke joinstyle="miter"/>
              <v:formulas>
                <v:f eqn="val #0"/>
              </v:formulas>
              <v:path arrowok="t" fillok="f" o:connecttype="none"/>
              <v:handles>
                <v:h position="#0,center"/>
              </v:handles>
              <o:lock v:ext="edit" shapetype="t"/>
            </v:shapetype>
            <v:rect id="_x0000_s1037" style="position:absolute;left:6823;top:258;width:2095;height:332">
              <v:textbox>
                <w:txbxContent>
                  <w:p w:rsidR="00DC1D4E" w:rsidRPr="00A524AD" w:rsidRDefault="00DC1D4E" w:rsidP="00DC1D4E">
                    <w:pPr>
                      <w:jc w:val="center"/>
                      <w:rPr>
                        <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
                        <w:sz w:val="16"/>
                        <w:szCs w:val="16"/>
                      </w:rPr>
                    </w:pPr>
                    <w:r>
                      <w:rPr>
                        <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
                        <w:sz w:val="16"/>
                        <w:szCs w:val="16"/>
                      </w:rPr>
                      <w:t>MSCAVSS</w:t>
                    </w:r>
                  </w:p>
                </w:txbxContent>
              </v:textbox>
            </v:rect>
            <v:shape id="_x0000_s1047" type="#_x0000_t34" style="position:absolute;left:4953;top:424;width:1870;height:917;rotation:180;flip:y" o:connectortype="elbow" adj=",254405,-84343">
              <v:stroke startarrow="block" endarrow="block"/>
            </v:shape>
            <v:rect id="_x0000_s1049" style="position:absolute;left:2859;top:1161;width:2094;height:358">
              <v:textbox>
                <w:txbxContent>
                  <w:p w:rsidR="00DC1D4E" w:rsidRPr="00A524AD" w:rsidRDefault="00DC1D4E" w:rsidP="00DC1D4E">
                    <w:pPr>
                      <w:jc w:val="center"/>
                      <w:rPr>
                        <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
                        <w:sz w:val="16"/>
                        <w:szCs w:val="16"/>
                      </w:rPr>
                    </w:pPr>
                    <w:r>
                      <w:rPr>
                        <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
                        <w:sz w:val="16"/>
                        <w:szCs w:val="16"/>
                      </w:rPr>
                      <w:t>App-V Agent</w:t>
                    </w:r>
                  </w:p>
                </w:txbxContent>
              </v:textbox>
            </v:rect>
            <v:shapetype id="_x0000_t34" coordsize="21600,21600" o:spt="34" o:oned="t" adj="10800" path="m,l@0,0@0,21600,21600,21600e" filled="f">
              <v:stroke joinstyle="miter"/>
              <v:formulas>
                <v:f eqn="val #0"/>
              </v:formulas>
              <v:path arrowok="t" fillok="f" o:connecttype="none"/>
              <v:handles>
                <v:h position="#0,center"/>
              </v:handles>
              <o:lock v:ext="edit" shapetype="t"/>
            </v:shapetype>
            <v:shape id="_x0000_s1035" type="#_x0000_t34" style="position:absolute;left:2956;top:291;width:1;height:495;rotation:180" o:connectortype="elbow" adj="-7776000,-486628,77954400">
              <v:stroke startarrow="block" endarrow="block"/>
            </v:shape>
            </v:group>
        </w:pict>
      </w:r>
    </w:p>
'@

$target.InsertXML($xml)
Write-Output "Inserted the rotation:180 connector shape (_x0000_s1035) and its shapetype."
